$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 699.0769
$ws.Range("I12").Value = 870.7143
$ws.Range("K12").Value = 870.7143
$ws.Range("M12").Value = -700.7143

$ws.Range("H43").Value = 5063.846
$ws.Range("I43").Value = 2574.8
$ws.Range("K43").Value = 2574.8
$ws.Range("M43").Value = -2505.8

$ws.Range("H62").Value = 39286.918
$ws.Range("I62").Value = 50875.65
$ws.Range("K62").Value = 50875.65
$ws.Range("M62").Value = -50251.65

$ws.Range("H65").Value = 39286.918
$ws.Range("I65").Value = 50875.65
$ws.Range("K65").Value = 254378.25
$ws.Range("M65").Value = -251258.25

$ws.Range("H70").Value = 183337630
$ws.Range("I70").Value = 50000876
$ws.Range("J70").Value = 250006000
$ws.Range("K70").Value = 150002628
$ws.Range("L70").Value = 750018000
$ws.Range("M70").Value = -150002358
$ws.Range("N70").Value = -750018540

$ws.Range("H73").Value = 183337630
$ws.Range("I73").Value = 50000876
$ws.Range("J73").Value = 250006000
$ws.Range("K73").Value = 150002628
$ws.Range("L73").Value = 750018000
$ws.Range("M73").Value = -150001692
$ws.Range("N73").Value = -750019872

$ws.Range("H86").Value = 3089.625
$ws.Range("J86").Value = 2633.8572
$ws.Range("L86").Value = 2633.8572
$ws.Range("N86").Value = -4879.8572

$ws.Range("H89").Value = 3089.625
$ws.Range("J89").Value = 2633.8572
$ws.Range("L89").Value = 13169.286
$ws.Range("N89").Value = -24401.286

$ws.Range("H106").Value = 5937.5
$ws.Range("I106").Value = 5539.1875
$ws.Range("J106").Value = 6999.6665
$ws.Range("K106").Value = 5539.1875
$ws.Range("L106").Value = 6999.6665
$ws.Range("M106").Value = -4908.1875
$ws.Range("N106").Value = -8261.666499999999

$ws.Range("H132").Value = 1457.4762
$ws.Range("I132").Value = 891.97144
$ws.Range("K132").Value = 2675.91432
$ws.Range("M132").Value = -145.9143199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6767.8945
$ws.Range("J32").Value = 6978.3335
$ws.Range("L32").Value = 6978.3335
$ws.Range("N32").Value = -7552.3335

$ws.Range("H45").Value = 2435.75
$ws.Range("I45").Value = 1926.8572
$ws.Range("J45").Value = 5998
$ws.Range("K45").Value = 1926.8572
$ws.Range("L45").Value = 5998
$ws.Range("M45").Value = -1549.8572
$ws.Range("N45").Value = -6752

$ws.Range("H74").Value = 1018.44446
$ws.Range("I74").Value = 1060.7059
$ws.Range("J74").Value = 300
$ws.Range("K74").Value = 1060.7059
$ws.Range("L74").Value = 300
$ws.Range("M74").Value = -186.7058999999999
$ws.Range("N74").Value = -2048

$ws.Range("H77").Value = 1018.44446
$ws.Range("I77").Value = 1060.7059
$ws.Range("J77").Value = 300
$ws.Range("K77").Value = 5303.5295
$ws.Range("L77").Value = 1500
$ws.Range("M77").Value = -935.5294999999996
$ws.Range("N77").Value = -10236

$ws.Range("H122").Value = 1951.36
$ws.Range("I122").Value = 1459.4375
$ws.Range("J122").Value = 2825.889
$ws.Range("K122").Value = 4378.3125
$ws.Range("L122").Value = 8477.667000000001
$ws.Range("M122").Value = -1928.3125
$ws.Range("N122").Value = -13377.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 340.88235
$ws.Range("J80").Value = 331.5
$ws.Range("L80").Value = 331.5
$ws.Range("N80").Value = -2327.5

$ws.Range("H83").Value = 340.88235
$ws.Range("J83").Value = 331.5
$ws.Range("L83").Value = 1657.5
$ws.Range("N83").Value = -11641.5

$ws.Range("H92").Value = 33332.668
$ws.Range("J92").Value = 33332.668
$ws.Range("L92").Value = 33332.668
$ws.Range("N92").Value = -38324.668

$ws.Range("H99").Value = 2667.04
$ws.Range("I99").Value = 1606.3334
$ws.Range("K99").Value = 1606.3334
$ws.Range("M99").Value = -108.3334

$ws.Range("H105").Value = 2006.6666
$ws.Range("I105").Value = 2006.6666
$ws.Range("K105").Value = 2006.6666
$ws.Range("M105").Value = -259.6666

$ws.Range("H122").Value = 77890
$ws.Range("J122").Value = 77890
$ws.Range("L122").Value = 77890
$ws.Range("N122").Value = -87690

$ws.Range("H134").Value = 1965.1305
$ws.Range("I134").Value = 2018.091
$ws.Range("K134").Value = 6054.272999999999
$ws.Range("M134").Value = -3519.272999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4109.3477
$ws.Range("I31").Value = 1786.7693
$ws.Range("K31").Value = 1786.7693
$ws.Range("M31").Value = -1491.7693

$ws.Range("H34").Value = 4109.3477
$ws.Range("I34").Value = 1786.7693
$ws.Range("K34").Value = 1786.7693
$ws.Range("M34").Value = -1584.7693

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 6499.5
$ws.Range("I87").Value = 6249.75
$ws.Range("K87").Value = 18749.25
$ws.Range("M87").Value = -17501.25

$ws.Range("H90").Value = 6499.5
$ws.Range("I90").Value = 6249.75
$ws.Range("K90").Value = 56247.75
$ws.Range("M90").Value = -50007.75

$ws.Range("H118").Value = 1907.591
$ws.Range("I118").Value = 1241.75
$ws.Range("J118").Value = 2055.5557
$ws.Range("K118").Value = 3725.25
$ws.Range("L118").Value = 6166.6671
$ws.Range("M118").Value = -2482.25
$ws.Range("N118").Value = -8652.667099999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7494.4443
$ws.Range("I80").Value = 4450
$ws.Range("K80").Value = 4450
$ws.Range("M80").Value = -3452

$ws.Range("H83").Value = 7494.4443
$ws.Range("I83").Value = 4450
$ws.Range("K83").Value = 22250
$ws.Range("M83").Value = -17258

$ws.Range("H113").Value = 4520.4443
$ws.Range("I113").Value = 2421.8
$ws.Range("J113").Value = 9290.091
$ws.Range("K113").Value = 2421.8
$ws.Range("L113").Value = 9290.091
$ws.Range("M113").Value = -251.8000000000002
$ws.Range("N113").Value = -13630.091

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4137.9375
$ws.Range("I7").Value = 3999.7693
$ws.Range("K7").Value = 3999.7693
$ws.Range("M7").Value = -3887.7693

$ws.Range("H40").Value = 5165.4614
$ws.Range("I40").Value = 4238.1904
$ws.Range("K40").Value = 4238.1904
$ws.Range("M40").Value = -4102.1904

$ws.Range("H61").Value = 3013.0967
$ws.Range("I61").Value = 1851.1
$ws.Range("J61").Value = 5125.8184
$ws.Range("K61").Value = 1851.1
$ws.Range("L61").Value = 5125.8184
$ws.Range("M61").Value = -1649.1
$ws.Range("N61").Value = -5529.8184

$ws.Range("H113").Value = 3013.0967
$ws.Range("I113").Value = 1851.1
$ws.Range("J113").Value = 5125.8184
$ws.Range("K113").Value = 1851.1
$ws.Range("L113").Value = 5125.8184
$ws.Range("M113").Value = 318.9000000000001
$ws.Range("N113").Value = -9465.8184

$ws.Range("H126").Value = 4137.9375
$ws.Range("I126").Value = 3999.7693
$ws.Range("K126").Value = 11999.3079
$ws.Range("M126").Value = -9529.3079

$ws.Range("H132").Value = 5942.8184
$ws.Range("I132").Value = 6688.75
$ws.Range("K132").Value = 20066.25
$ws.Range("M132").Value = -17536.25

$ws.Range("H136").Value = 15754.139
$ws.Range("I136").Value = 1278.0435
$ws.Range("J136").Value = 41365.69
$ws.Range("K136").Value = 3834.1305
$ws.Range("L136").Value = 124097.07
$ws.Range("M136").Value = -1284.1305
$ws.Range("N136").Value = -129197.07

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 19830
$ws.Range("J63").Value = 19830
$ws.Range("L63").Value = 19830
$ws.Range("N63").Value = -21078

$ws.Range("H66").Value = 19830
$ws.Range("J66").Value = 19830
$ws.Range("L66").Value = 59490
$ws.Range("N66").Value = -65730

$ws.Range("H107").Value = 1564.8148
$ws.Range("I107").Value = 1531.25
$ws.Range("J107").Value = 1833.3334
$ws.Range("K107").Value = 4593.75
$ws.Range("L107").Value = 5500.0002
$ws.Range("M107").Value = -2673.75
$ws.Range("N107").Value = -9340.0002

$ws.Range("H113").Value = 1070.5
$ws.Range("J113").Value = 1924.25
$ws.Range("L113").Value = 5772.75
$ws.Range("N113").Value = -10112.75

$ws.Range("H136").Value = 3245.1562
$ws.Range("I136").Value = 3012.652
$ws.Range("J136").Value = 3839.3333
$ws.Range("K136").Value = 9037.956
$ws.Range("L136").Value = 11517.9999
$ws.Range("M136").Value = -6487.956
$ws.Range("N136").Value = -16617.9999
